$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Return_with_prediction (G), return_pct_change (H), and mean_return_pct_change (I, row 2 only)
# values produced by a re-run of the auto-recurrence model for the semester input-250 S&P500 comparison.
$updates = @(
    @{ Cell = "G2"; Value = 0.06224091725287884 }
    @{ Cell = "H2"; Value = -3.877415028843533 }
    @{ Cell = "I2"; Value = 222.2636742179514 }
    @{ Cell = "G3"; Value = 0.0757114981186691 }
    @{ Cell = "H3"; Value = 10.68111834055966 }
    @{ Cell = "G4"; Value = -0.03325172432826021 }
    @{ Cell = "H4"; Value = 25.93122174440089 }
    @{ Cell = "G5"; Value = -0.02159632838261233 }
    @{ Cell = "H5"; Value = 20.22586434019264 }
    @{ Cell = "G6"; Value = -0.1083859720864143 }
    @{ Cell = "H6"; Value = -2.211191303456858 }
    @{ Cell = "G7"; Value = -0.09113813621269763 }
    @{ Cell = "H7"; Value = 0.2550385426666495 }
    @{ Cell = "G8"; Value = -0.3624880555821954 }
    @{ Cell = "H8"; Value = 1.203711997819769 }
    @{ Cell = "G9"; Value = -0.3860728810120887 }
    @{ Cell = "H9"; Value = 1.029976736516805 }
    @{ Cell = "G10"; Value = 0.01996457991918173 }
    @{ Cell = "H10"; Value = -1.121380968284853 }
    @{ Cell = "G11"; Value = 0.02726826057822341 }
    @{ Cell = "H11"; Value = 20.15373487306747 }
    @{ Cell = "G12"; Value = 0.2123672471662505 }
    @{ Cell = "H12"; Value = -4.23017244300395 }
    @{ Cell = "G13"; Value = 0.2280844804021571 }
    @{ Cell = "H13"; Value = 1.277794012793083 }
    @{ Cell = "G14"; Value = -0.04974051409339369 }
    @{ Cell = "H14"; Value = -18.13417113575828 }
    @{ Cell = "G15"; Value = -0.04849114657065651 }
    @{ Cell = "H15"; Value = -1.666481773803567 }
    @{ Cell = "G16"; Value = 0.2117606893631132 }
    @{ Cell = "H16"; Value = -0.3823194313013588 }
    @{ Cell = "G17"; Value = 0.2223837316021025 }
    @{ Cell = "H17"; Value = 0.8323664223632846 }
    @{ Cell = "G18"; Value = 0.07487194779065426 }
    @{ Cell = "H18"; Value = 2.535882084479091 }
    @{ Cell = "G19"; Value = 0.07362470182070974 }
    @{ Cell = "H19"; Value = -2.274165330987812 }
    @{ Cell = "G20"; Value = -0.08492938609975433 }
    @{ Cell = "H20"; Value = -13.26640148202525 }
    @{ Cell = "G21"; Value = -0.07978073040038806 }
    @{ Cell = "H21"; Value = 7.841275984479799 }
    @{ Cell = "G22"; Value = 0.07723222572854922 }
    @{ Cell = "H22"; Value = 5.075146431700891 }
    @{ Cell = "G23"; Value = 0.07415785074534212 }
    @{ Cell = "H23"; Value = 8.526887381654774 }
    @{ Cell = "G24"; Value = 0.06054380588139363 }
    @{ Cell = "H24"; Value = -9.109940187452208 }
    @{ Cell = "G25"; Value = 0.06998761401182149 }
    @{ Cell = "H25"; Value = 27.76112648675191 }
    @{ Cell = "G26"; Value = 0.1139806638290622 }
    @{ Cell = "H26"; Value = -4.498042591722923 }
    @{ Cell = "G27"; Value = 0.1182927581574682 }
    @{ Cell = "H27"; Value = 3.892998179843158 }
    @{ Cell = "G28"; Value = 0.1311033994688534 }
    @{ Cell = "H28"; Value = 1.431313116346654 }
    @{ Cell = "G29"; Value = 0.1387189121307623 }
    @{ Cell = "H29"; Value = -8.035977054699064 }
    @{ Cell = "G30"; Value = 0.08656343611097742 }
    @{ Cell = "H30"; Value = 2.676077874762662 }
    @{ Cell = "G31"; Value = 0.0911034209695494 }
    @{ Cell = "H31"; Value = 11.52601494228961 }
    @{ Cell = "G32"; Value = 0.05539078954472102 }
    @{ Cell = "H32"; Value = 3.806027674249955 }
    @{ Cell = "G33"; Value = 0.05680614817648612 }
    @{ Cell = "H33"; Value = 2.82975194107236 }
    @{ Cell = "G34"; Value = 0.01512511749045012 }
    @{ Cell = "H34"; Value = -12.85991375785732 }
    @{ Cell = "G35"; Value = 0.01891392993710746 }
    @{ Cell = "H35"; Value = 11.9121883351389 }
    @{ Cell = "G36"; Value = -0.02965581390376378 }
    @{ Cell = "H36"; Value = -2.100529141449747 }
    @{ Cell = "G37"; Value = -0.02463890391516297 }
    @{ Cell = "H37"; Value = 25.92895087022834 }
    @{ Cell = "G38"; Value = 0.08404500693154533 }
    @{ Cell = "H38"; Value = 7.360862930669635 }
    @{ Cell = "G39"; Value = 0.07327633957935985 }
    @{ Cell = "H39"; Value = -5.748001572591319 }
    @{ Cell = "G40"; Value = 0.06868150554948413 }
    @{ Cell = "H40"; Value = 3.734345392282935 }
    @{ Cell = "G41"; Value = 0.07475740134099189 }
    @{ Cell = "H41"; Value = 14.97098860551124 }
    @{ Cell = "G42"; Value = 0.07732335203782546 }
    @{ Cell = "H42"; Value = -0.6015340163498276 }
    @{ Cell = "G43"; Value = 0.09408946455477447 }
    @{ Cell = "H43"; Value = 17.37049927329974 }
    @{ Cell = "G44"; Value = 0.08944399182462591 }
    @{ Cell = "H44"; Value = 1.356063815299923 }
    @{ Cell = "G45"; Value = 0.08414275987911836 }
    @{ Cell = "H45"; Value = -6.904739636240109 }
    @{ Cell = "G46"; Value = 0.002899601702761789 }
    @{ Cell = "H46"; Value = 205.9710178971165 }
    @{ Cell = "G47"; Value = 0.005738710367670439 }
    @{ Cell = "H47"; Value = 12112.35586066376 }
    @{ Cell = "G48"; Value = -0.1038631224974459 }
    @{ Cell = "H48"; Value = -8.072386027512978 }
    @{ Cell = "G49"; Value = -0.1057895793805522 }
    @{ Cell = "H49"; Value = 3.45569896909119 }
    @{ Cell = "G50"; Value = 0.1636766468100906 }
    @{ Cell = "H50"; Value = -4.001567453605435 }
    @{ Cell = "G51"; Value = 0.171575586710011 }
    @{ Cell = "H51"; Value = 1.024947409649907 }
    @{ Cell = "G52"; Value = 0.06081105279337013 }
    @{ Cell = "H52"; Value = -14.30462544291136 }
    @{ Cell = "G53"; Value = 0.06910153735101483 }
    @{ Cell = "H53"; Value = 7.440483259517194 }
    @{ Cell = "G54"; Value = -0.1339093449465831 }
    @{ Cell = "H54"; Value = -4.778053874472053 }
    @{ Cell = "G55"; Value = -0.1147886650882252 }
    @{ Cell = "H55"; Value = 1.444565668204366 }
    @{ Cell = "G56"; Value = 0.1875715094621933 }
    @{ Cell = "H56"; Value = -1.288977422520355 }
    @{ Cell = "G57"; Value = 0.2037098029156213 }
    @{ Cell = "H57"; Value = 2.413941084349778 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
